$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values for column A, rows 1-27 ---
$values = @(
    "cod",
    7898994908722,
    7898994908739,
    7898994908746,
    7898994908753,
    7898994908715,
    8777,
    5064,
    5071,
    5088,
    8760,
    "TESTESHELFLIFE",
    7896579902028,
    5293,
    5286,
    5255,
    5279,
    5262,
    3342,
    3341,
    5354,
    5361,
    5378,
    5309,
    5316,
    5323,
    5378
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Remove any leftover formatting/content below row 27 (none expected, but be safe)
$ws.Rows.Item(28).Resize(1000).ClearContents()

# --- Formatting ---

# Row 1 header "cod": bold, size 16, Consolas, vertical center, row height 21
$header = $ws.Range("A1")
$header.Font.Name = "Consolas"
$header.Font.Bold = $true
$header.Font.Size = 16
$header.VerticalAlignment = -4108   # xlVAlignCenter
$ws.Rows.Item(1).RowHeight = 21

# Rows 2-25: plain number format "0", default font, no border, no alignment override
$body = $ws.Range("A2:A25")
$body.Borders.LineStyle = -4142      # xlLineStyleNone
$body.NumberFormat = "0"

# Row 26: red OpenSansRegular size 10 font, General number format, no border
$r26 = $ws.Range("A26")
$r26.Font.Name = "OpenSansRegular"
$r26.Font.Size = 10
$r26.Font.Color = 255

# Row 27: number format "0", red Calibri font, thin border all sides, right aligned
$r27 = $ws.Range("A27")
$r27.NumberFormat = "0"
$r27.Font.Color = 255
$r27.HorizontalAlignment = -4152   # xlHAlignRight
$r27.Borders.LineStyle = 1         # xlContinuous
$r27.Borders.Weight = 2            # xlThin

$ws.Range("H10:H11").Select()
